$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.385.29"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "1.709.02"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'224.32"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("D6").Value = "'0.5337"
$ws.Range("E6").Value = "  -2.18%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.2675"
$ws.Range("E8").Value = "  -2.55%  "
$ws.Range("D9").Value = "'0.06605"
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("D10").Value = "'20.95"
$ws.Range("E10").Value = "  -4.29%  "
$ws.Range("D11").Value = "'0.07629"
$ws.Range("E11").Value = "  -2.04%  "
$ws.Range("D12").Value = "'4.553"
$ws.Range("D13").Value = "1.731.51"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").Value = "1.945.27"
$ws.Range("E14").Value = "  -1.44%  "
$ws.Range("D15").Value = "'0.5773"
$ws.Range("E15").Value = "  -3.60%  "
$ws.Range("D16").Value = "0.0₅8181"
$ws.Range("E16").Value = "  -2.88%  "
$ws.Range("D17").Value = "'67.76"
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("D18").Value = "27.348.96"
$ws.Range("E18").Value = "  -1.27%  "
$ws.Range("D19").Value = "'216.69"
$ws.Range("E19").Value = "  -5.07%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "'4.670"
$ws.Range("E21").Value = "  -3.37%  "
$ws.Range("E22").Value = "  -4.21%  "
$ws.Range("D23").Value = "'5.965"
$ws.Range("E23").Value = "  -4.22%  "
$ws.Range("D24").Value = "'1.004"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'142.21"
$ws.Range("E25").Value = "  -4.06%  "
$ws.Range("D26").Value = "'1.740"
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("E27").Value = "  -3.03%  "
$ws.Range("D28").Value = "'7.262"
$ws.Range("E28").Value = "  -2.88%  "
$ws.Range("E29").Value = "  -5.10%  "
$ws.Range("D30").Value = "'0.05395"
$ws.Range("E30").Value = "  -5.46%  "
$ws.Range("D32").Value = "'3.493"
$ws.Range("E32").Value = "  -5.66%  "
$ws.Range("D33").Value = "'3.427"
$ws.Range("E33").Value = "  -2.66%  "
$ws.Range("D34").Value = "'1.641"
$ws.Range("E34").Value = "  -2.86%  "
$ws.Range("D35").Value = "'2.874"
$ws.Range("E35").Value = "  +0.70%  "
$ws.Range("D36").Value = "'0.9495"
$ws.Range("E36").Value = "  -2.75%  "
$ws.Range("D37").Value = "'2.415"
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("D38").Value = "'0.5860"
$ws.Range("E38").Value = "  -2.20%  "
$ws.Range("D39").Value = "'0.01636"
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("D40").Value = "'5.852"
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("D41").Value = "1.045.11"
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "'0.8407"
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("D44").Value = "'100.94"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("D45").Value = "1.851.96"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("E46").Value = "  +2.10%  "
$ws.Range("D47").Value = "'57.95"
$ws.Range("E47").Value = "  -2.89%  "
$ws.Range("D48").Value = "'0.4515"
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").Value = "'8.082"
$ws.Range("E50").Value = "  -2.95%  "
$ws.Range("E51").Value = "  -1.86%  "
